$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.969.48"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "1.677.19"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'214.98"
$ws.Range("E5").Value = "  -0.63%  "
$ws.Range("E6").Value = "  +1.45%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  -0.26%  "
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("D10").Value = "'20.32"
$ws.Range("E11").Value = "  -0.23%  "
$ws.Range("D12").Value = "1.914.15"
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("D13").Value = "1.675.38"
$ws.Range("E13").Value = "  -0.07%  "
$ws.Range("E14").Value = "  -0.22%  "
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").Value = "26.987.09"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("D18").Value = "'237.13"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").Value = "'8.06"
$ws.Range("E19").Value = "  +3.93%  "
$ws.Range("E20").Value = "  -0.85%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("E22").Value = "  -0.76%  "
$ws.Range("E23").Value = "  -1.05%  "
$ws.Range("D24").Value = "'2.19"
$ws.Range("E24").Value = "  -1.83%  "
$ws.Range("D25").Value = "'145.55"
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("E26").Value = "  +1.43%  "
$ws.Range("D27").Value = "'16.07"
$ws.Range("E27").Value = "  +0.57%  "
$ws.Range("E28").Value = "  -1.56%  "
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("E30").Value = "  -0.38%  "
$ws.Range("E31").Value = "  -0.88%  "
$ws.Range("E32").Value = "  +0.12%  "
$ws.Range("D33").Value = "1.485.21"
$ws.Range("E33").Value = "  +0.54%  "
$ws.Range("E34").Value = "  +0.61%  "
$ws.Range("D35").Value = "'1.68"
$ws.Range("E35").Value = "  +3.69%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("D37").Value = "'0.583"
$ws.Range("E37").Value = "  +1.27%  "
$ws.Range("D38").Value = "'0.0175"
$ws.Range("E38").Value = "  +2.79%  "
$ws.Range("D39").Value = "'0.899"
$ws.Range("E39").Value = "  -0.13%  "
$ws.Range("E40").Value = "  -3.50%  "
$ws.Range("E41").Value = "  +0.41%  "
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("D43").Value = "'2.31"
$ws.Range("E43").Value = "  +1.26%  "
$ws.Range("D44").Value = "'67.49"
$ws.Range("D45").Value = "1.819.25"
$ws.Range("E45").Value = "  -0.34%  "
$ws.Range("E46").Value = "  -0.20%  "
$ws.Range("D47").Value = "'90.54"
$ws.Range("E47").Value = "  +0.22%  "
$ws.Range("E48").Value = "  +13.98%  "
$ws.Range("E49").Value = "  -0.63%  "
$ws.Range("E50").Value = "  +1.30%  "
$ws.Range("E51").Value = "  +0.40%  "
